$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the indicator title (trailing period after "5.6.1.1" removed)
$ws.Range("A1").Value = "5.6.1.1 15-49 жаштардын азыркы замандагы контрацепция методу менен кабардар болгон күйөөгө тийген жана күйөөгө тийбеген сексуалдык активдүү аялдардын үлүшү"
$ws.Range("B1").Value = "5.6.1.1 Доля замужних женщин и сексуально активных не замужних женщин в возрасте 15-49 лет, которые были осведомлены о соврменном методе контрацепции"
$ws.Range("C1").Value = "5.6.1.1 Proportion of married women and sexually active single women aged 15-49 years who were aware of the modern method of contraception"

# Rows 6-7: replace short "urban/rural" labels with the fuller phrasing
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the active selection
$ws.Range("A6:C7").Select()
